# The underlying edit swaps the two theme parts in the package:
#   ppt/theme/theme1.xml (was "Office Theme")  <->  ppt/theme/theme2.xml (was "Integral")
# so that, after the edit, the presentation's primary/applied theme (referenced by
# the slide master and by the presentation itself) is the plain "Office Theme"
# palette, while the theme that used to back the slide master ("Integral") ends up
# parked on the file that is no longer wired to the slides.
#
# This PowerPoint COM host only exposes one writable theme palette for a deck -
# the one attached to the SlideMaster / Presentation (ppt/theme/theme2.xml in this
# file) - via Master.ColorScheme / NotesMaster.ColorScheme (they both resolve to the
# same underlying document theme). There is no COM surface to address the second,
# unused theme part directly. So we reproduce the reachable, visible half of the
# swap: push the "Office Theme" color palette (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) into the presentation's applied color scheme, in place of the
# "Integral" colors that were there before.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$cs = $master.ColorScheme

# Index -> (scheme slot, target "Office Theme" RGB)
# COM RGB longs are packed as R + G*256 + B*65536 (i.e. 0x00BBGGRR).
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
